{"js": "// Homework0.docx edit:\n//   1. \"Due: 1/26/17\" -> \"Due: 1/25/17\", with the text split into two runs\n//      (\"Due: 1/25\" + \"/17\") separated by a moved \"_GoBack\" bookmark (this\n//      is exactly what Word does when you retype the day digits and save:\n//      the last-edit-position bookmark lands between the edited text and\n//      the untouched tail).\n//   2. The old \"_GoBack\" bookmark (previously sitting at the very end of\n//      the document, after the MAC-address paragraph) is removed, since\n//      Word keeps only a single \"_GoBack\" bookmark per document and moves\n//      it to the most recent edit location.\n\nconst body = context.document.body;\n\n// --- Step 1: fix the typo \"26\" -> \"25\" in the \"Due:\" line ------------------\nconst dateHits = body.search(\"Due: 1/26/17\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\n\nif (dateHits.items.length > 0) {\n  const dateRange = dateHits.items[0];\n  dateRange.insertText(\"Due: 1/25/17\", \"Replace\");\n  await context.sync();\n}\n\n// --- Step 2: drop the bookmark that used to sit at the document's end -----\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Step 3: re-insert \"_GoBack\" between \"Due: 1/25\" and \"/17\" ------------\nconst dueHits = body.search(\"Due: 1/25\", { matchCase: true });\ndueHits.load(\"items\");\nawait context.sync();\n\nif (dueHits.items.length > 0) {\n  const insertionPoint = dueHits.items[0].getRange(\"End\");\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Homework0.docx edit:\n#   1. \"Due: 1/26/17\" -> \"Due: 1/25/17\", with the text split into two runs\n#      (\"Due: 1/25\" + \"/17\") separated by a moved \"_GoBack\" bookmark (this\n#      is exactly what Word does when you retype the day digits and save:\n#      the last-edit-position bookmark lands between the edited text and\n#      the untouched tail).\n#   2. The old \"_GoBack\" bookmark (previously sitting at the very end of\n#      the document, after the MAC-address paragraph) is removed, since\n#      Word keeps only a single \"_GoBack\" bookmark per document and moves\n#      it to the most recent edit location.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: drop the bookmark that used to sit at the document's end -----\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Step 2: fix the typo \"26\" -> \"25\" in the \"Due:\" line ------------------\n$dateRng = $d.Content\n$dateRng.Find.Execute(\"Due: 1/26/17\", $false, $false, $false, $false, $false, $true, 1, $false, \"Due: 1/25/17\", 2)\n\n# --- Step 3: re-insert \"_GoBack\" between \"Due: 1/25\" and \"/17\" ------------\n$insRng = $d.Content\n$insRng.Find.Execute(\"Due: 1/25\")\n$insRng.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $insRng)\n"}
